$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ensure cells stay formatted as text so values are written verbatim (matching the source data format).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "293.51"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-5.74%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "40.55"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.11%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.029"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.60%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07329"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-3.58%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.543"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-8.18%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9285"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.60%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1166"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.84%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1742"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-4.35%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.04334"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "4.73%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08644"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-4.76%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1055"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.02%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001267"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.19%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.006027"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2.86%"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.01%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.287"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.82%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.82%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.972"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "5.15%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "4.37%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-3.31%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.03940"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.01%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001261"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.12%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.003658"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-10.00%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-5.18%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "22.49%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02314"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-4.98%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05074"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-2.06%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.006204"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "87.92%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007857"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.84%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1288"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.05%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007357"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.06%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007269"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-14.24%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3197"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-4.84%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006278"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-4.72%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.01%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.03329"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-87.85%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.01%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.01%"
